# BUON ANNO SI VOLAAAAAA
# Insert a new slide "Assunzioni fatte 4" right after "Assunzioni fatte 3"
# (position 8 of 11 in the final deck), pushing "Considerazione su processi
# lavorativi", "Esecuzione degli scostamenti" and "Possibili ipotesi agli
# scostamenti" down by one slot. No other slide content changes.

$p = $ppt.ActivePresentation

# Insert the new slide at index 8 using the same "Titolo e contenuto"
# (Title and Content) layout used by the other "Assunzioni fatte" slides.
$newSlide = $p.Slides.Add(8, $ppt.ppLayoutText)

# Match the Italian placeholder names used throughout the rest of the deck.
$title = $newSlide.Shapes.Item(1)
$title.Name = "Titolo 1"

$body = $newSlide.Shapes.Item(2)
$body.Name = "Segnaposto contenuto 2"

# Title text.
$title.TextFrame.TextRange.Text = "Assunzioni fatte 4"
$title.TextFrame.TextRange.LanguageID = "it-IT"

# Body text: two paragraphs, the second split across two runs.
$bodyText = $body.TextFrame.TextRange
$para1 = "Se produco meno di quello che vendo si pensa che ci sia un magazzino col prodotto gi$([char]0x00E0) disponibile e quindi come quantit$([char]0x00E0) si prende il numero di pezzi prodotti."
$para2Run1 = "Se invece produco di pi$([char]0x00F9) di quello che vendo si assume che ci sia un surplus della produzione per fare magazzino / prevenzione errori e quindi per coerenza di scelta  si utilizza come quantit$([char]0x00E0) per calcolare i costi il numero di "
$para2Run2 = "pezzi prodotti."

$bodyText.Text = $para1 + "`r" + $para2Run1
$bodyText.InsertAfter($para2Run2) | Out-Null
$bodyText.LanguageID = "it-IT"
